$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Texts of tenders that are no longer present in the freshly scraped data.
$removedTexts = @(
    "Objet : Gardiennage et surveillance des Bâtiments Administratifs relevant de la Direction Régionale de l’Artisanat et de l’Economie Sociale Souss Massa et ses entités en lot unique.",
    "Objet : LA CONCEPTION ARCHITECTURALE ET LE SUIVI DES TRAVAUX DE CONSTRUCTION DU BLOC FONCIER DE L’ANCFCC DE GUERCIF"
)

# Collect the data rows (everything below the header row) that should be kept.
$keptA = New-Object System.Collections.ArrayList
$keptB = New-Object System.Collections.ArrayList
for ($r = 2; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Text
    $b = $ws.Cells.Item($r, 2).Text
    if ($removedTexts -notcontains $a) {
        [void]$keptA.Add($a)
        [void]$keptB.Add($b)
    }
}

# Newly scraped tenders go to the top of the data, most recent first.
$newA = @(
    "Objet : Acquisition des disjoncteurs 72,5 kV pour la Division Exploitation Transport Marrakech",
    "Objet : Peinture des lignes HTB du réseau de Division Exploitation Transport Agadir Lot n° 1 : Peinture des lignes 60 & 225 kV Lot n° 2 : Peinture des lignes 400 kV"
)
$newB = @("N/A", "N/A")

$finalA = $newA + $keptA
$finalB = $newB + $keptB

for ($i = 0; $i -lt $finalA.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $finalA[$i]
    $ws.Cells.Item($r, 2).Value = $finalB[$i]
}

# If the new data set is shorter than the old one, clear out the now-stale trailing rows.
$newLastRow = $finalA.Count + 1
if ($newLastRow -lt $lastRow) {
    $clearRange = $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($lastRow, 2))
    $clearRange.Clear()
}
